$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- begroting (sheet1) ---

# Row 19: apply the built-in "Neutral" cell style to B19 (Part 5 entry)
$ws1.Range("B19").Style = "Neutral"

# Row 26: relabel A26 as "week number" and move its old label ("hours spent") to B26
$ws1.Range("A26").Value = "week number"
$ws1.Range("B26").Value = "hours spent"

# New rows 38-39: extra weekly samples
$ws1.Range("A38").Value = 15
$ws1.Range("B38").Value = 34
$ws1.Range("A39").Value = 16

# New row 41: column header for the budget total below
$ws1.Range("C41").Value = "budget"

# B42 / C42 keep their existing formulas; they recalc automatically once
# B38 participates in SUM(B27:B40)

# --- Sheet1 (sheet2) ---
$ws2.Range("B6").Value = 5

# --- view/selection state ---
$ws2.Range("B9").Select() | Out-Null
$ws1.Activate()
$ws1.Range("B36").Select() | Out-Null
